$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update total right count (B12: 45 -> 75)
$ws.Range("B12").Value = 75

# Update correct/total marks label (E12: "37/84" -> "75/140")
$ws.Range("E12").Value = "75/140"
